# GMT synchronized with eVITTA server
#
# Slide 1, "TextBox 4" shape contains a list of CSV filenames, one per
# paragraph. The paragraph for "df_proteomic_scale.csv" was previously
# split across four separate runs ("df" / "_" / "proteomic_" /
# "scale.csv") - this collapses it back into a single run while keeping
# the original (first-run) character formatting.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$targetText = "df_proteomic_scale.csv"

# Locate the paragraph holding the (possibly multi-run) target text.
$paraIndex = 0
$paraCount = $tr.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $candidate = $tr.Paragraphs($i, 1)
    $candidateText = $candidate.Text.TrimEnd("`r", "`n")
    if ($candidateText -eq $targetText) {
        $paraIndex = $i
        break
    }
}

if ($paraIndex -eq 0) {
    throw "Could not locate paragraph containing '$targetText'"
}

$para = $tr.Paragraphs($paraIndex, 1)

# Re-assigning the exact same text is treated as a no-op (the existing
# run split survives untouched), so first swap in an unrelated
# placeholder string - this collapses the paragraph down to a single
# run carrying the first original run's formatting - then assign the
# real text back onto that single run.
$para.Text = "PLACEHOLDER"

$para = $sh.TextFrame.TextRange.Paragraphs($paraIndex, 1)
$para.Text = $targetText
